# Atualização automática de preços de eletricidade
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = 45910
$ws.Range("B2").Value = 89.94
$ws.Range("C2").Value = 83.09
$ws.Range("D2").Value = 83.09
$ws.Range("E2").Value = 81.09999999999999
$ws.Range("F2").Value = 72.98999999999999
$ws.Range("G2").Value = 66
$ws.Range("H2").Value = 77.5
$ws.Range("I2").Value = 84.26000000000001
$ws.Range("J2").Value = 75.98999999999999
$ws.Range("K2").Value = 45.1
$ws.Range("L2").Value = 8
$ws.Range("M2").Value = 0
$ws.Range("N2").Value = 0
$ws.Range("O2").Value = 0
$ws.Range("P2").Value = 0
$ws.Range("Q2").Value = 0
$ws.Range("R2").Value = 0
$ws.Range("S2").Value = 4.9
$ws.Range("T2").Value = 23.09
$ws.Range("U2").Value = 62.17
$ws.Range("V2").Value = 96
$ws.Range("W2").Value = 105.01
$ws.Range("X2").Value = 96.26000000000001
$ws.Range("Y2").Value = 83.09
$ws.Range("Z2").Value = 51.57
$ws.Range("AB2").Value = 95.09
$ws.Range("AD2").Value = 100.5
$ws.Range("AF2").Value = 89.68000000000001
$ws.Range("AG2").Value = "9h-18h"
